$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NASA Floater Platform CAPEX")

# ------------------------------------------------------------------
# Preserve formatting of the existing "totals" block (rows 12-14) by
# copying it down to its new location (rows 14-16) BEFORE overwriting
# row 12, since row 14 is both a copy source and a copy destination.
# ------------------------------------------------------------------
$ws.Range("A14:D14").Copy($ws.Range("A16"))
$ws.Range("A13:C13").Copy($ws.Range("A15"))
$ws.Range("A12:C12").Copy($ws.Range("A14"))

# Row 13 becomes a blank spacer row again (fully cleared - no leftover cells).
$ws.Range("A13:D13").Clear()

# D14 held the inflation note that moved to D16; clear the original.
$ws.Range("D14").ClearContents()

# Z column formatting continues down through the new rows.
$ws.Range("Z10").Copy($ws.Range("Z11"))
$ws.Range("Z10").Copy($ws.Range("Z12"))

# --- New cell F5: ratio used for the Aqua Ventus / CPI comparison ---
$ws.Range("F5").Formula = "=22200000/15000"

# --- B6 gets the same (currently empty) numeric style as B5 ---
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("B6").HorizontalAlignment = $ws.Range("B5").HorizontalAlignment
$ws.Range("B6").VerticalAlignment = $ws.Range("B5").VerticalAlignment

# --- Row 19: Comparison cases heading ---
$ws.Range("A19").Value = "Comparison cases"

# --- Row 20: Aqua Ventus comparison, 2018 USD ---
$ws.Range("A20").Value = "AQUA VENTUS PLATFORM CAPEX ($/KW)"
$ws.Range("B20").Value = 690
$ws.Range("C20").Value = "2018 USD"

# --- Row 21: Aqua Ventus comparison, 2021 USD (inflated) ---
$ws.Range("A21").Value = "AQUA VENTUS PLATFORM CAPEX ($/KW)"
$ws.Range("B21").Value = 728
$ws.Range("C21").Value = "2021 USD"
$ws.Range("D21").Value = "Inflated using BLS CPI from 1/1/2018 to 1/1/2021(https://www.bls.gov/data/inflation_calculator.htm)"

# --- Row 12 repurposed: TMD sensitivity ($) = TMD CAPEX * TMD scale ---
$ws.Range("A12").Value = "TMD sensitivity ($)"
$ws.Range("A12").Font.Bold = $false
$ws.Range("B12").Formula = "=B10*B11"
$ws.Range("B12").Font.Bold = $false
$ws.Range("C12").Clear()

# --- New CPI ratio cells ---
$ws.Range("E14").Value = "CPI ratio"
$ws.Range("F14").Formula = "=369/350"

# --- New row 11: TMD scale factor (B11 takes B10's numeric style) ---
$ws.Range("A11").Value = "TMD sclae"
$ws.Range("B10").Copy($ws.Range("B11"))
$ws.Range("B11").Value = 1

# --- Row 14: TOTAL PLATFORM CAPEX ($) now sums B9+B12, C14 left blank ---
$ws.Range("A14").Value = "TOTAL PLATFORM CAPEX ($)"
$ws.Range("B14").Formula = "=B9+B12"
$ws.Range("C14").ClearContents()

# --- Row 15: TOTAL PLATFORM CAPEX ($/kW), 2018 USD ---
$ws.Range("A15").Value = "TOTAL PLATFORM CAPEX ($/kW)"
$ws.Range("B15").Formula = "=B14/B2"
$ws.Range("C15").Value = "2018 USD"

# --- Row 16: TOTAL PLATFORM CAPEX ($/kW), 2021 USD (inflated) ---
$ws.Range("A16").Value = "TOTAL PLATFORM CAPEX ($/kW)"
$ws.Range("B16").Formula = "=B15*F14"
$ws.Range("C16").Value = "2021 USD"
$ws.Range("D16").Value = "Inflated using BLS CPI from 1/1/2018 to 1/1/2021(https://www.bls.gov/data/inflation_calculator.htm)"

# --- Row 17: total platform CAPEX scaled by turbine count ---
$ws.Range("B17").Formula = "=B16*15000"

# --- Selection / view updates ---
$ws.Range("G6").Select()

$ws2 = $wb.Worksheets.Item("Sheet1")
$ws2.Range("P20").Select()
